# Auto-generated edit script applying the Goblin_Profits market-price refresh diff.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets, matching the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 46667.24
$ws.Range("I62").Value = 105561
$ws.Range("J62").Value = 10424.923
$ws.Range("K62").Value = 105561
$ws.Range("L62").Value = 10424.923
$ws.Range("M62").Value = -104937
$ws.Range("N62").Value = -11672.923
$ws.Range("H64").Value = 7967.6875
$ws.Range("I64").Value = 3524
$ws.Range("K64").Value = 3524
$ws.Range("M64").Value = -3276
$ws.Range("H65").Value = 46667.24
$ws.Range("I65").Value = 105561
$ws.Range("J65").Value = 10424.923
$ws.Range("K65").Value = 527805
$ws.Range("L65").Value = 52124.61500000001
$ws.Range("M65").Value = -524685
$ws.Range("N65").Value = -58364.61500000001
$ws.Range("H67").Value = 7967.6875
$ws.Range("I67").Value = 3524
$ws.Range("K67").Value = 3524
$ws.Range("M67").Value = -2666
$ws.Range("H70").Value = 5561018
$ws.Range("I70").Value = 25004250
$ws.Range("J70").Value = 5808.857
$ws.Range("K70").Value = 75012750
$ws.Range("L70").Value = 17426.571
$ws.Range("M70").Value = -75012480
$ws.Range("N70").Value = -17966.571
$ws.Range("H73").Value = 5561018
$ws.Range("I73").Value = 25004250
$ws.Range("J73").Value = 5808.857
$ws.Range("K73").Value = 75012750
$ws.Range("L73").Value = 17426.571
$ws.Range("M73").Value = -75011814
$ws.Range("N73").Value = -19298.571
$ws.Range("H80").Value = 380.16666
$ws.Range("J80").Value = 532.1429000000001
$ws.Range("L80").Value = 1596.4287
$ws.Range("N80").Value = -3592.4287
$ws.Range("H83").Value = 380.16666
$ws.Range("J83").Value = 532.1429000000001
$ws.Range("L83").Value = 4789.2861
$ws.Range("N83").Value = -14773.2861
$ws.Range("H137").Value = 14075.875
$ws.Range("J137").Value = 1751.5
$ws.Range("L137").Value = 5254.5
$ws.Range("N137").Value = -10354.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2704.375
$ws.Range("I2").Value = 1652.2222
$ws.Range("K2").Value = 1652.2222
$ws.Range("M2").Value = -1539.2222
$ws.Range("H116").Value = 2704.375
$ws.Range("I116").Value = 1652.2222
$ws.Range("K116").Value = 1652.2222
$ws.Range("M116").Value = 641.7778000000001
$ws.Range("H122").Value = 1776.4773
$ws.Range("I122").Value = 1544.2433
$ws.Range("K122").Value = 4632.7299
$ws.Range("M122").Value = -2182.7299

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2704.375
$ws.Range("I3").Value = 1652.2222
$ws.Range("K3").Value = 1652.2222
$ws.Range("M3").Value = -1538.2222
$ws.Range("H103").Value = 18984
$ws.Range("J103").Value = 18984
$ws.Range("L103").Value = 18984
$ws.Range("N103").Value = -21328

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4000
$ws.Range("I16").Value = 4000
$ws.Range("K16").Value = 4000
$ws.Range("M16").Value = -3713
$ws.Range("H31").Value = 6951.1333
$ws.Range("I31").Value = 2085.5557
$ws.Range("K31").Value = 2085.5557
$ws.Range("M31").Value = -1790.5557
$ws.Range("H34").Value = 6951.1333
$ws.Range("I34").Value = 2085.5557
$ws.Range("K34").Value = 2085.5557
$ws.Range("M34").Value = -1883.5557
$ws.Range("H59").Value = 37000
$ws.Range("I59").Value = 9000
$ws.Range("J59").Value = 65000
$ws.Range("K59").Value = 9000
$ws.Range("L59").Value = 65000
$ws.Range("M59").Value = -7855
$ws.Range("N59").Value = -67290
$ws.Range("H62").Value = 10401
$ws.Range("I62").Value = 11000
$ws.Range("K62").Value = 11000
$ws.Range("M62").Value = -10376
$ws.Range("H65").Value = 10401
$ws.Range("I65").Value = 11000
$ws.Range("K65").Value = 55000
$ws.Range("M65").Value = -51880
$ws.Range("H68").Value = 58333.332
$ws.Range("J68").Value = 67500
$ws.Range("L68").Value = 67500
$ws.Range("N68").Value = -68998
$ws.Range("H71").Value = 58333.332
$ws.Range("J71").Value = 67500
$ws.Range("L71").Value = 202500
$ws.Range("N71").Value = -209988
$ws.Range("H99").Value = 2791.3333
$ws.Range("I99").Value = 2843.5
$ws.Range("J99").Value = 2687
$ws.Range("K99").Value = 2843.5
$ws.Range("L99").Value = 2687
$ws.Range("M99").Value = -1345.5
$ws.Range("N99").Value = -5683
$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 4000
$ws.Range("K113").Value = 4000
$ws.Range("M113").Value = -1830
$ws.Range("H126").Value = 2791.3333
$ws.Range("I126").Value = 2843.5
$ws.Range("J126").Value = 2687
$ws.Range("K126").Value = 8530.5
$ws.Range("L126").Value = 8061
$ws.Range("M126").Value = -6060.5
$ws.Range("N126").Value = -13001
$ws.Range("H134").Value = 2293
$ws.Range("I134").Value = 2081.0952
$ws.Range("K134").Value = 6243.285600000001
$ws.Range("M134").Value = -3708.285600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 124
$ws.Range("I2").Value = 130
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 780
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -667
$ws.Range("N2").Value = -826

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 484.63635
$ws.Range("I107").Value = 464
$ws.Range("J107").Value = 539.6667
$ws.Range("K107").Value = 464
$ws.Range("L107").Value = 539.6667
$ws.Range("M107").Value = 1456
$ws.Range("N107").Value = -4379.6667
$ws.Range("H136").Value = 62333
$ws.Range("J136").Value = 62333
$ws.Range("L136").Value = 186999
$ws.Range("N136").Value = -192099

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5062.1816
$ws.Range("I7").Value = 4860.5
$ws.Range("J7").Value = 5600
$ws.Range("K7").Value = 4860.5
$ws.Range("L7").Value = 5600
$ws.Range("M7").Value = -4748.5
$ws.Range("N7").Value = -5824
$ws.Range("H16").Value = 494.45
$ws.Range("I16").Value = 507.8421
$ws.Range("J16").Value = 240
$ws.Range("K16").Value = 507.8421
$ws.Range("L16").Value = 240
$ws.Range("M16").Value = -337.8421
$ws.Range("N16").Value = -580
$ws.Range("H126").Value = 5062.1816
$ws.Range("I126").Value = 4860.5
$ws.Range("J126").Value = 5600
$ws.Range("K126").Value = 14581.5
$ws.Range("L126").Value = 16800
$ws.Range("M126").Value = -12111.5
$ws.Range("N126").Value = -21740
$ws.Range("H132").Value = 4657.769
$ws.Range("I132").Value = 4356.0625
$ws.Range("J132").Value = 5140.5
$ws.Range("K132").Value = 13068.1875
$ws.Range("L132").Value = 15421.5
$ws.Range("M132").Value = -10538.1875
$ws.Range("N132").Value = -20481.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 425
$ws.Range("J4").Value = 425
$ws.Range("L4").Value = 425
$ws.Range("N4").Value = -651
$ws.Range("H100").Value = 658.75
$ws.Range("I100").Value = 492.5
$ws.Range("J100").Value = 825
$ws.Range("K100").Value = 985
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -444
$ws.Range("N100").Value = -2732
$ws.Range("H126").Value = 2288
$ws.Range("I126").Value = 2285.8
$ws.Range("J126").Value = 2299
$ws.Range("K126").Value = 6857.400000000001
$ws.Range("L126").Value = 6897
$ws.Range("M126").Value = -4387.400000000001
$ws.Range("N126").Value = -11837
$ws.Range("H136").Value = 2547.6365
$ws.Range("I136").Value = 1838.25
$ws.Range("J136").Value = 4439.3335
$ws.Range("K136").Value = 5514.75
$ws.Range("L136").Value = 13318.0005
$ws.Range("M136").Value = -2964.75
$ws.Range("N136").Value = -18418.0005
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
